# Apply the edit described by the diff:
# - Fill D4:D16 with "Pass" (same shared string already used in D2:D3)
# - Scroll the view so row 7 is at the top
# - Move the active selection to D16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Pass" values in column D for rows 4 through 16
$ws.Range("D4:D16").Value = "Pass"

# Scroll the window so that A7 becomes the top-left visible cell
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1

# Update the active selection to D16
$ws.Range("D16").Select()
